$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text values in the source data (e.g. "56.836.96"
# with thousands separators). Several of the new prices look like plain
# numbers (e.g. "492.36"), and a bare `.Value =` assignment would let Excel
# auto-convert those to real numbers. Force text with NumberFormat "@"
# before assigning, then ClearFormats so we don't leave a stray explicit
# "Text" cell style behind (the source cells carry no style at all).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "56.748.25"
$ws.Range("E2").Value = "  +1.89%  "

Set-TextValue $ws.Range("D3") "2.497.43"
$ws.Range("E3").Value = "  -1.02%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue $ws.Range("D5") "492.36"
$ws.Range("E5").Value = "  +0.20%  "

Set-TextValue $ws.Range("D6") "152.46"
$ws.Range("E6").Value = "  +7.53%  "

$ws.Range("E7").Value = "  -0.05%  "

Set-TextValue $ws.Range("D8") "0.514"
$ws.Range("E8").Value = "  -0.22%  "

Set-TextValue $ws.Range("D9") "2.512.47"
$ws.Range("E9").Value = "  -0.34%  "

Set-TextValue $ws.Range("D10") "5.74"
$ws.Range("E10").Value = "  +3.21%  "

Set-TextValue $ws.Range("D11") "0.0989"
$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("E13").Value = "  +0.82%  "

Set-TextValue $ws.Range("D14") "2.931.66"
$ws.Range("E14").Value = "  -0.24%  "

Set-TextValue $ws.Range("D15") "56.811.91"
$ws.Range("E15").Value = "  +1.85%  "

Set-TextValue $ws.Range("D16") "21.28"
$ws.Range("E16").Value = "  +1.41%  "

Set-TextValue $ws.Range("D17") "0.0000137"
$ws.Range("E17").Value = "  -1.48%  "

Set-TextValue $ws.Range("D18") "2.509.58"
$ws.Range("E18").Value = "  -0.31%  "

Set-TextValue $ws.Range("D19") "4.57"
$ws.Range("E19").Value = "  +3.13%  "

Set-TextValue $ws.Range("D20") "10.33"
$ws.Range("E20").Value = "  +2.21%  "

Set-TextValue $ws.Range("D21") "322.62"
$ws.Range("E21").Value = "  -0.45%  "

$ws.Range("E22").Value = "  +0.01%  "

Set-TextValue $ws.Range("D23") "5.91"
$ws.Range("E23").Value = "  +1.84%  "

Set-TextValue $ws.Range("D24") "58.68"
$ws.Range("E24").Value = "  +0.45%  "

Set-TextValue $ws.Range("D25") "0.412"
$ws.Range("E25").Value = "  -0.57%  "

Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  -0.57%  "

Set-TextValue $ws.Range("D27") "0.162"
$ws.Range("E27").Value = "  -5.45%  "

Set-TextValue $ws.Range("D28") "2.598.72"
$ws.Range("E28").Value = "  -0.66%  "

Set-TextValue $ws.Range("D29") "7.59"
$ws.Range("E29").Value = "  +0.98%  "

Set-TextValue $ws.Range("D30") "0.0₃0807"
$ws.Range("E30").Value = "  +0.41%  "

Set-TextValue $ws.Range("D31") "0.998"
$ws.Range("E31").Value = "  -0.02%  "

Set-TextValue $ws.Range("D32") "151.08"
$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("E34").Value = "  +0.70%  "

Set-TextValue $ws.Range("D35") "5.28"
$ws.Range("E35").Value = "  +0.52%  "

$ws.Range("E36").Value = "  +2.35%  "

Set-TextValue $ws.Range("D37") "3.79"
$ws.Range("E37").Value = "  +1.28%  "

Set-TextValue $ws.Range("D38") "0.873"
$ws.Range("E38").Value = "  -0.96%  "

Set-TextValue $ws.Range("D39") "1.40"
$ws.Range("E39").Value = "  +4.37%  "

Set-TextValue $ws.Range("D40") "34.24"
$ws.Range("E40").Value = "  -0.76%  "

$ws.Range("E41").Value = "  +1.86%  "

Set-TextValue $ws.Range("D42") "0.0565"
$ws.Range("E42").Value = "  +1.33%  "

Set-TextValue $ws.Range("D43") "0.617"
$ws.Range("E43").Value = "  -0.05%  "

Set-TextValue $ws.Range("D44") "0.995"
$ws.Range("E44").Value = "  -0.03%  "

# Rows 45 and 46 swap: Bittensor <-> RenderToken (each keeping its own
# refreshed price/volume figures) change places in the ranking.
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D45") "4.89"
$ws.Range("E45").Value = "  +2.00%  "

$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D46") "268.67"
$ws.Range("E46").Value = "  +3.39%  "

Set-TextValue $ws.Range("D47") "0.0935"
$ws.Range("E47").Value = "  +1.96%  "

$ws.Range("E48").Value = "  +1.03%  "

Set-TextValue $ws.Range("D49") "10.24"
$ws.Range("E49").Value = "  +0.91%  "

Set-TextValue $ws.Range("D50") "17.92"
$ws.Range("E50").Value = "  +1.53%  "

Set-TextValue $ws.Range("D51") "1.903.59"
$ws.Range("E51").Value = "  -5.62%  "
